{"js": "// Find the paragraph that carries the \"{{ Vorname }} ... also took part in\n// workshops with the following companies:\" placeholder (currently split\n// across six runs), and locate the \"_Hlk96973631\" bookmark that currently\n// starts on the following paragraph (the one with the\n// \"{%- if workshops|length > 2 %}\" Jinja tag) and ends much further down,\n// right after \"{% endif %}\".\nconst results = context.document.body.search(\n  \"also took part in workshops with the following companies:\",\n  { matchCase: false }\n);\nresults.load(\"items\");\nawait context.sync();\n\nconst vornamePara = results.items[0].paragraphs.getFirst();\n\n// Collapse the six separate runs (\"{{\", \" \", \"Vorname\", \" \", \"}}\", \" also\n// took part...\") into a single run carrying the new Jinja condition, and\n// drop the paragraph's justified alignment (the replacement paragraph in the\n// target has no <w:jc/>).\nconst wholeRange = vornamePara.getRange(\"Whole\");\nwholeRange.insertText(\n  \"{%- if workshops|length > 0 %}{{ Vorname }} also took part in workshops with the following companies: {% endif %}\",\n  Word.InsertLocation.replace\n);\nvornamePara.alignment = Word.Alignment.left;\nawait context.sync();\n\n// Move just the bookmark's *start* up onto this paragraph, keeping its\n// existing end position untouched: build a range that spans from the start\n// of the rewritten paragraph through to the bookmark's current end, delete\n// the old bookmark, then re-create it over that full span so the\n// <w:bookmarkStart> lands here while <w:bookmarkEnd> stays where it was.\nconst doc = context.document;\nconst existingBookmarkRange = doc.getBookmarkRange(\"_Hlk96973631\");\nconst startOfVornamePara = vornamePara.getRange(\"Start\");\nconst newBookmarkSpan = startOfVornamePara.expandTo(existingBookmarkRange);\nawait context.sync();\n\ndoc.deleteBookmark(\"_Hlk96973631\");\nnewBookmarkSpan.insertBookmark(\"_Hlk96973631\");\nawait context.sync();\n", "ps1": "# Locate the paragraph carrying the \"{{ Vorname }} ... also took part in\n# workshops with the following companies:\" placeholder (currently split\n# across six runs: \"{{\", \" \", \"Vorname\", \" \", \"}}\", \" also took part...\").\n$d = $word.ActiveDocument\n\n$findRange = $d.Content\n$found = $findRange.Find.Execute(\"also took part in workshops with the following companies:\")\nif (-not $found) {\n    throw \"Could not find the workshops placeholder paragraph\"\n}\n\n$vornamePara = $findRange.Paragraphs(1)\n$paraRange = $vornamePara.Range\n\n# Replace the paragraph's own text (excluding its trailing paragraph mark) so\n# the six runs collapse into a single run; the new run inherits the first\n# run's formatting (font/size), matching the target markup.\n$textRange = $d.Range($paraRange.Start, $paraRange.End - 1)\n$textRange.Text = \"{%- if workshops|length > 0 %}{{ Vorname }} also took part in workshops with the following companies: {% endif %}\"\n\n# Drop the paragraph's justified alignment - the rewritten paragraph in the\n# target has no <w:jc/> element.\n$vornamePara.Format.Alignment = 0\n\n# Move just the \"_Hlk96973631\" bookmark's start up onto this paragraph while\n# leaving its existing end position (far below, right after \"{% endif %}\")\n# untouched: capture the bookmark's current end, delete the bookmark, then\n# re-add it spanning from the start of the rewritten paragraph through to\n# that same end point.\n$bookmark = $d.Bookmarks(\"_Hlk96973631\")\n$bookmarkEnd = $bookmark.End\n$paraStart = $vornamePara.Range.Start\n\n$d.Bookmarks(\"_Hlk96973631\").Delete()\n$newBookmarkRange = $d.Range($paraStart, $bookmarkEnd)\n$d.Bookmarks.Add(\"_Hlk96973631\", $newBookmarkRange)\n"}
